$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "You need to make [separate] branches" - merge the "separate" run (and
#    its surrounding gramStart/gramEnd proofErr markers + trailing space run)
#    back into a single run with the preceding text. Re-"typing" the exact
#    same text over the whole span collapses it into one run and drops the
#    now-unnecessary proofErr markers.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "You have a local repository but you can" + [char]8217 + "t just shove work in there yet! You need to make separate ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You have a local repository but you can" + [char]8217 + "t just shove work in there yet! You need to make separate ",
    2)

# ---------------------------------------------------------------------------
# 2) Fix "-prune" -> "--prune" typo, and move the "_GoBack" bookmark from the
#    end of the document into the middle of "from" ("fro" + bookmark + "m").
# ---------------------------------------------------------------------------

# Remove the existing "_GoBack" bookmark (currently sitting alone in an empty
# paragraph near the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate "git fetch --all -prune" (the second, advanced-use mention).
$rng = $d.Content
$found = $rng.Find.Execute("git fetch --all -prune")
$gitStart = $rng.Start

# Insert the missing "-" right before the lone "-" that precedes "prune",
# turning "-prune" into "--prune".
$insPoint = $d.Range($gitStart + 16, $gitStart + 16)
$insPoint.InsertBefore("-")

# Split "git fetch --all --prune" into three runs:
#   "git fetch --all " | "-" | "-prune"
# Adding (and immediately removing) a temporary bookmark at a point inside a
# run forces the run to split there, while deleting the bookmark again
# leaves the split in place without leaving any bookmark behind.
$split1 = $d.Range($gitStart + 16, $gitStart + 16)
$d.Bookmarks.Add("_TmpSplitA", $split1) | Out-Null
$d.Bookmarks.Item("_TmpSplitA").Delete()

$split2 = $d.Range($gitStart + 17, $gitStart + 17)
$d.Bookmarks.Add("_TmpSplitB", $split2) | Out-Null
$d.Bookmarks.Item("_TmpSplitB").Delete()

# Re-find the first mention, "which grabs all changes from the upstream", to
# drop the "_GoBack" bookmark right in the middle of "from" ("fro" | "m").
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(" which grabs all changes fro")
$goBackPoint = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null
